# Apply the "Fruta / hortaliza, semanal" weekly-refresh edit to the Cilantro subset sheet.
# The source feed rotates one weekly record out of the window and one new record in:
#   - a brand-new record (Primera+Segunda) is inserted at the front of the series (rows 60-61)
#   - every later record shifts down by one slot (rows 62-167 take on the date that used to
#     sit two rows above them), carrying along the handful of non-date fields (Volumen, Origen)
#     that differ from this market/category's usual values
#   - the oldest record that fell out of the rolling window is re-appended at the end (rows 168-169)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fecha (column D) updates for rows 60-167 ---
$ws.Range("D60").Value = 44614
$ws.Range("D61").Value = 44614
$ws.Range("D62").Value = 44306
$ws.Range("D63").Value = 44306
$ws.Range("D64").Value = 44222
$ws.Range("D65").Value = 44222
$ws.Range("D66").Value = 44383
$ws.Range("D67").Value = 44383
$ws.Range("D68").Value = 44469
$ws.Range("D69").Value = 44469
$ws.Range("D70").Value = 44278
$ws.Range("D71").Value = 44278
$ws.Range("D72").Value = 44435
$ws.Range("D73").Value = 44435
$ws.Range("D74").Value = 44292
$ws.Range("D75").Value = 44292
$ws.Range("D76").Value = 44467
$ws.Range("D77").Value = 44467
$ws.Range("D78").Value = 44203
$ws.Range("D79").Value = 44203
$ws.Range("D80").Value = 44341
$ws.Range("D81").Value = 44341
$ws.Range("D82").Value = 44245
$ws.Range("D83").Value = 44245
$ws.Range("D84").Value = 44442
$ws.Range("D85").Value = 44442
$ws.Range("D86").Value = 44350
$ws.Range("D87").Value = 44350
$ws.Range("D88").Value = 44574
$ws.Range("D89").Value = 44574
$ws.Range("D90").Value = 44398
$ws.Range("D91").Value = 44398
$ws.Range("D92").Value = 44509
$ws.Range("D93").Value = 44509
$ws.Range("D94").Value = 44433
$ws.Range("D95").Value = 44433
$ws.Range("D96").Value = 44316
$ws.Range("D97").Value = 44316
$ws.Range("D98").Value = 44582
$ws.Range("D99").Value = 44582
$ws.Range("D100").Value = 44159
$ws.Range("D101").Value = 44159
$ws.Range("D102").Value = 44386
$ws.Range("D103").Value = 44386
$ws.Range("D104").Value = 44322
$ws.Range("D105").Value = 44322
$ws.Range("D106").Value = 44320
$ws.Range("D107").Value = 44320
$ws.Range("D108").Value = 44475
$ws.Range("D109").Value = 44475
$ws.Range("D110").Value = 44204
$ws.Range("D111").Value = 44204
$ws.Range("D112").Value = 44358
$ws.Range("D113").Value = 44358
$ws.Range("D114").Value = 44313
$ws.Range("D115").Value = 44313
$ws.Range("D116").Value = 44460
$ws.Range("D117").Value = 44460
$ws.Range("D118").Value = 44194
$ws.Range("D119").Value = 44194
$ws.Range("D120").Value = 44217
$ws.Range("D121").Value = 44217
$ws.Range("D122").Value = 44166
$ws.Range("D123").Value = 44166
$ws.Range("D124").Value = 44579
$ws.Range("D125").Value = 44579
$ws.Range("D126").Value = 44237
$ws.Range("D127").Value = 44237
$ws.Range("D128").Value = 44609
$ws.Range("D129").Value = 44609
$ws.Range("D130").Value = 44252
$ws.Range("D131").Value = 44252
$ws.Range("D132").Value = 44271
$ws.Range("D133").Value = 44271
$ws.Range("D134").Value = 44420
$ws.Range("D135").Value = 44420
$ws.Range("D136").Value = 44336
$ws.Range("D137").Value = 44336
$ws.Range("D138").Value = 44231
$ws.Range("D139").Value = 44231
$ws.Range("D140").Value = 44565
$ws.Range("D141").Value = 44565
$ws.Range("D142").Value = 44334
$ws.Range("D143").Value = 44334
$ws.Range("D144").Value = 44280
$ws.Range("D145").Value = 44280
$ws.Range("D146").Value = 44362
$ws.Range("D147").Value = 44362
$ws.Range("D148").Value = 44365
$ws.Range("D149").Value = 44365
$ws.Range("D150").Value = 44567
$ws.Range("D151").Value = 44567
$ws.Range("D152").Value = 44553
$ws.Range("D153").Value = 44553
$ws.Range("D154").Value = 44490
$ws.Range("D155").Value = 44490
$ws.Range("D156").Value = 44462
$ws.Range("D157").Value = 44462
$ws.Range("D158").Value = 44264
$ws.Range("D159").Value = 44264
$ws.Range("D160").Value = 44330
$ws.Range("D161").Value = 44330
$ws.Range("D162").Value = 44257
$ws.Range("D163").Value = 44257
$ws.Range("D164").Value = 44299
$ws.Range("D165").Value = 44299
$ws.Range("D166").Value = 44285
$ws.Range("D167").Value = 44285

# --- Volumen (column J) follow-along updates for rows 70-73 ---
$ws.Range("J70").Value = 200
$ws.Range("J71").Value = 100
$ws.Range("J72").Value = 400
$ws.Range("J73").Value = 200

# --- Origen (column O) follow-along updates ---
$ws.Range("O94").Value = 'Región de Ñuble'
$ws.Range("O95").Value = 'Región de Ñuble'
$ws.Range("O96").Value = 'Región Metropolitana'
$ws.Range("O97").Value = 'Región Metropolitana'
$ws.Range("O104").Value = 'Región de Ñuble'
$ws.Range("O105").Value = 'Región de Ñuble'
$ws.Range("O106").Value = 'Región Metropolitana'
$ws.Range("O107").Value = 'Región Metropolitana'
$ws.Range("O122").Value = 'Región de Ñuble'
$ws.Range("O123").Value = 'Región de Ñuble'
$ws.Range("O124").Value = 'Región Metropolitana'
$ws.Range("O125").Value = 'Región Metropolitana'

# --- Append new rows 168-169 (the record that rolled off the end of the 60-167 window) ---
$ws.Range("A168").Value = 11
$ws.Range("B168").Value = 'Vega Monumental Concepción'
$ws.Range("C168").Value = 'Bíobío'
$ws.Range("D168").Value = 44595
$ws.Range("D168").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E168").Value = 8
$ws.Range("F168").Value = 100112040
$ws.Range("G168").Value = 'Cilantro'
$ws.Range("H168").Value = 'Sin especificar'
$ws.Range("I168").Value = 'Primera'
$ws.Range("J168").Value = 200
$ws.Range("K168").Value = 600
$ws.Range("L168").Value = 700
$ws.Range("M168").Value = 650
$ws.Range("N168").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O168").Value = 'Región de Ñuble'
$ws.Range("P168").Value = 650
$ws.Range("Q168").Value = 1
$ws.Range("R168").Value = 'Hortaliza'
$ws.Range("A169").Value = 11
$ws.Range("B169").Value = 'Vega Monumental Concepción'
$ws.Range("C169").Value = 'Bíobío'
$ws.Range("D169").Value = 44595
$ws.Range("D169").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E169").Value = 8
$ws.Range("F169").Value = 100112040
$ws.Range("G169").Value = 'Cilantro'
$ws.Range("H169").Value = 'Sin especificar'
$ws.Range("I169").Value = 'Segunda'
$ws.Range("J169").Value = 100
$ws.Range("K169").Value = 500
$ws.Range("L169").Value = 500
$ws.Range("M169").Value = 500
$ws.Range("N169").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O169").Value = 'Región de Ñuble'
$ws.Range("P169").Value = 500
$ws.Range("Q169").Value = 1
$ws.Range("R169").Value = 'Hortaliza'
